# "Verify Valid login test"
# The Valid_Login sheet gains a new "FirstName" column (inserted before the
# old "Type" column, which shifts right along with the "Vaild" value), and
# the sample login values are swapped out for a fresh valid-login test row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valid_Login")

# Insert a new column C (pushes old Type/Vaild column from C to D)
$ws.Columns.Item(3).Insert()

# Fill in the new test data for the updated "valid login" row
$ws.Range("A2").Value = "Login1"
$ws.Range("B2").Value = "Pwd@abcd1"
$ws.Range("C2").Value = "Fname1"
$ws.Range("C1").Value = "FirstName"

# The inserted column picked up formatting from its left neighbour; reset it
# back to the default so it matches a plain, unstyled column.
$ws.Range("C1:C2").Style = "Normal"

# Resize column B to fit its new (longer) contents
$ws.Columns.Item(2).EntireColumn.AutoFit()

# Update the sheet's active cell/selection
$ws.Range("F11").Select() | Out-Null
